$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '91.412.56'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('E2').Value = '  +2.51%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.124.87'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('E3').Value = '  +0.92%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.51'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +2.66%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '624.75'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  +0.43%  '
$ws.Range('E7').Value = '  +27.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.376'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  +0.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.121.53'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  +0.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.772'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  +23.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.193'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  +6.75%  '
$ws.Range('E13').Value = '  +4.60%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.04'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  +8.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.54'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  +4.39%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.171.56'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.698.63'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  +0.93%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.91'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  +15.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.104.48'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  +0.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000217'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  +2.53%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.15'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  +5.21%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '443.44'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  +4.70%  '
$ws.Range('E23').Value = '  +6.64%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.19'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  +4.88%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '6.36'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  +14.33%  '
$ws.Range('B26').Value = 'Aptos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.63'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  +6.19%  '
$ws.Range('B27').Value = 'Litecoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '88.71'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  +7.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.284.84'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  +1.52%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.167'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  -2.38%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.35'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  +15.06%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '530.43'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  +3.99%  '
$ws.Range('B33').Value = 'Binance-PegBSC-USD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.898'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  -16.20%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.75'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  +1.80%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.17'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  +6.25%  '
$ws.Range('E36').Value = '  +11.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '24.29'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  +8.76%  '
$ws.Range('E38').Value = '  +4.12%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.87'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  +3.82%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0885'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  +27.32%  '
$ws.Range('B41').Value = 'WhiteBITCoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '22.29'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('E42').Value = '  +21.37%  '
$ws.Range('E43').Value = '  -0.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.401'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  +10.17%  '
$ws.Range('B45').Value = 'Stacks'
$ws.Range('C45').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.94'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  +6.48%  '
$ws.Range('B46').Value = 'USDe'
$ws.Range('C46').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.00'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '149.21'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  +2.48%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '44.36'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  +2.39%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.32'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  +8.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.34'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  +9.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '168.53'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  +4.61%  '
